# Applies the cryptos-list price/volume refresh described by the commit:
# "Updated cryptos list on Tue Apr  9 16:50:31 UTC 2024 with GitHub Actions"
#
# Columns on Sheet1: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Both D and E hold plain text in the source workbook (note values like
# "68.773.63" and "4.066.63" are not valid numbers - they are thousands-dotted
# display strings), so every write below is a text write. Column E is always
# safe (it contains "%" and padding spaces so Excel can never mistake it for
# a number). Some column D values, though, are themselves valid numeric
# literals (e.g. "6.62", "0.999"); assigning those directly through
# Range.Value would make Excel auto-convert the cell to a Number. To keep
# them as text (matching the original inlineStr cell), we prefix with an
# apostrophe (Excel's literal "treat as text" marker) and then reset
# Style to "Normal" so the quote-prefix indicator doesn't leave a stray
# style applied to the cell (the source cells carry no style at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D) text or $null if unchanged, new Volume(1h)
# (E) text or $null if unchanged, and whether the Price text must be forced
# to stay text (i.e. it parses as a plain number).
$updates = @(
    @(2, "68.773.63", "  -4.31%  ", $false),
    @(3, "3.505.97", "  -4.13%  ", $false),
    @(4, $null, "  -0.03%  ", $false),
    @(5, "576.61", "  -2.60%  ", $true),
    @(6, "173.80", "  -4.10%  ", $true),
    @(7, $null, "  -0.61%  ", $false),
    @(8, "3.493.99", "  -4.39%  ", $false),
    @(9, $null, "  +0.04%  ", $false),
    @(10, "0.188", "  -7.39%  ", $true),
    @(11, "6.62", "  +6.78%  ", $true),
    @(12, "0.600", "  -1.82%  ", $true),
    @(13, "47.10", "  -5.71%  ", $true),
    @(14, "0.0000274", "  -4.45%  ", $true),
    @(15, "679.66", "  -0.47%  ", $true),
    @(16, "4.066.63", "  -4.13%  ", $false),
    @(17, "8.86", "  -2.38%  ", $true),
    @(18, "68.822.59", "  -4.36%  ", $false),
    @(19, "3.509.19", "  -4.24%  ", $false),
    @(20, $null, "  -1.46%  ", $false),
    @(21, "17.47", "  -4.29%  ", $true),
    @(22, "11.17", "  -4.17%  ", $true),
    @(23, "0.902", "  -4.40%  ", $true),
    @(24, "16.29", "  -9.48%  ", $true),
    @(25, "97.36", "  -6.02%  ", $true),
    @(26, "3.83", "  -5.26%  ", $true),
    @(27, "0.999", "  -0.12%  ", $true),
    @(28, $null, "  -6.74%  ", $false),
    @(29, "9.35", "  -8.51%  ", $true),
    @(30, "32.98", "  -6.82%  ", $true),
    @(31, "8.77", "  -4.86%  ", $true),
    @(32, $null, "  -9.32%  ", $false),
    @(33, $null, "  -6.07%  ", $false),
    @(34, "7.23", "  -1.56%  ", $true),
    @(35, "562.03", "  -3.13%  ", $true),
    @(36, "3.61", "  -14.28%  ", $true),
    @(37, "10.83", "  -4.55%  ", $true),
    @(38, $null, "  -3.71%  ", $false),
    @(39, "57.13", "  -4.07%  ", $true),
    @(41, $null, "  -4.58%  ", $false),
    @(42, "0.0439", "  -6.01%  ", $true),
    @(43, "3.449.26", "  -7.58%  ", $false),
    @(44, "0.334", "  -3.81%  ", $true),
    @(45, "33.35", "  -6.69%  ", $true),
    @(46, "0.0₃0700", "  -8.26%  ", $false),
    @(47, "2.89", "  +2.32%  ", $true),
    @(48, $null, "  -7.96%  ", $false),
    @(49, $null, "  -0.60%  ", $false),
    @(50, "134.50", "  +1.31%  ", $true),
    @(51, $null, "  -0.88%  ", $false)
)

foreach ($u in $updates) {
    $row = $u[0]
    $priceVal = $u[1]
    $volVal = $u[2]
    $forceText = $u[3]

    if ($null -ne $priceVal) {
        $priceAddr = "D" + $row
        if ($forceText) {
            $ws.Range($priceAddr).Value = "'" + $priceVal
            $ws.Range($priceAddr).Style = "Normal"
        } else {
            $ws.Range($priceAddr).Value = $priceVal
        }
    }

    if ($null -ne $volVal) {
        $volAddr = "E" + $row
        $ws.Range($volAddr).Value = $volVal
    }
}
